# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values, column G, recalculated from the regenerated
# source data (replacing the old Strike# based values) keyed by row number.
$kValues = @{
    2 = 4
    3 = 1
    4 = 0
    5 = 1
    6 = 0
    7 = 1
    8 = 0
    9 = 1
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 2
    34 = 0
    35 = 2
    36 = 0
    37 = 4
    38 = 0
    39 = 3
    40 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 1
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 2
    54 = 1
    55 = 2
    56 = 1
    57 = 0
    58 = 3
    59 = 1
    60 = 1
    61 = 1
    62 = 1
    63 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
